# Generate Report for Handoff
# Regenerates the localization-status report: the two dependent .png rows
# are dropped, the single .md source file is replaced by two newly
# handed-off .md files (new GUID-based names), and all related
# handoff/handback metadata (hashes, timestamps, hyperlinks) is refreshed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# New file identities used across the whole workbook
# ---------------------------------------------------------------------
$md1 = "27191251-2634-40f4-a04b-7b7a3a1b317e.md"
$md2 = "eb5b6d7f-1b9b-4f04-88e3-fae0653135b2.md"

$xlf1zh = "27191251-2634-40f4-a04b-7b7a3a1b317e.4e9d883f10139ba1bc316ffd214f987f2fd835fa.zh-cn.xlf"
$xlf2zh = "eb5b6d7f-1b9b-4f04-88e3-fae0653135b2.c5f9a66ec1bfd70a8e49530b4826e8162b806280.zh-cn.xlf"
$xlf1de = "27191251-2634-40f4-a04b-7b7a3a1b317e.4e9d883f10139ba1bc316ffd214f987f2fd835fa.de-de.xlf"
$xlf2de = "eb5b6d7f-1b9b-4f04-88e3-fae0653135b2.c5f9a66ec1bfd70a8e49530b4826e8162b806280.de-de.xlf"

$handoffZh = "2016-03-09 03:40:13"
$handoffDe = "2016-03-09 03:40:25"
$epoch = "0001-01-01 00:00:00"

$commitSrc = "2cbc9ca12d17b7c537ca3b03c4d1f44ce3e0e377"
$commitZh  = "a3c24ace3f63c66f6bf529304e0ae55d4dbf1b32"
$commitDe  = "ae111d72fe686f85743e6695fc1bb36c78306c4d"

$urlMd1      = "https://github.com/OpenLocalizationTest/oltest/blob/$commitSrc/e2e/$md1"
$urlMd2      = "https://github.com/OpenLocalizationTest/oltest/blob/$commitSrc/e2e/$md2"
$urlConfig   = "https://github.com/OpenLocalizationTest/oltest/blob/$commitSrc/.localization-config"
$urlXlf1Zh   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitZh/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$xlf1zh"
$urlXlf2Zh   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitZh/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$xlf2zh"
$urlXlf1De   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitDe/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$xlf1de"
$urlXlf2De   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitDe/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$xlf2de"

# =======================================================================
# Sheet 1: Overview
# =======================================================================
$ws1 = $wb.Worksheets.Item(1)

# Drop the row belonging to the second dependent picture (no longer present)
$ws1.Rows.Item(5).Delete()

$ws1.Range("A2").Value = $md1
$ws1.Range("A3").Value = $md2
$ws1.Range("A4").Value = ".localization-config"
$ws1.Range("B4").Value = "Not to be localized"
$ws1.Range("C4").Value = "Not to be localized"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), $urlMd1, "", "", $md1)
$ws1.Hyperlinks.Add($ws1.Range("A3"), $urlMd2, "", "", $md2)
$ws1.Hyperlinks.Add($ws1.Range("A4"), $urlConfig, "", "", ".localization-config")

# =======================================================================
# Sheet 2: zh-cn
# =======================================================================
$ws2 = $wb.Worksheets.Item(2)

$ws2.Rows.Item(5).Delete()

$ws2.Range("A2").Value = $md1
$ws2.Range("C2").Value = $xlf1zh
$ws2.Range("D2").Value = $handoffZh
$ws2.Range("H2").Value = "Include"

$ws2.Range("A3").Value = $md2
$ws2.Range("C3").Value = $xlf2zh
$ws2.Range("D3").Value = $handoffZh
$ws2.Range("H3").Value = "Include"

$ws2.Range("A4").Value = ".localization-config"
$ws2.Range("B4").Value = "Not to be localized"
$ws2.Range("C4").ClearContents()
$ws2.Range("D4").Value = $epoch
$ws2.Range("H4").Value = "Ignored"

$ws2.Range("I2").ClearContents()
$ws2.Range("I3").ClearContents()

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $urlMd1, "", "", $md1)
$ws2.Hyperlinks.Add($ws2.Range("C2"), $urlXlf1Zh, "", "", $xlf1zh)
$ws2.Hyperlinks.Add($ws2.Range("A3"), $urlMd2, "", "", $md2)
$ws2.Hyperlinks.Add($ws2.Range("C3"), $urlXlf2Zh, "", "", $xlf2zh)
$ws2.Hyperlinks.Add($ws2.Range("A4"), $urlConfig, "", "", ".localization-config")

# =======================================================================
# Sheet 3: de-de
# =======================================================================
$ws3 = $wb.Worksheets.Item(3)

$ws3.Rows.Item(5).Delete()

$ws3.Range("A2").Value = $md1
$ws3.Range("C2").Value = $xlf1de
$ws3.Range("D2").Value = $handoffDe
$ws3.Range("H2").Value = "Include"

$ws3.Range("A3").Value = $md2
$ws3.Range("C3").Value = $xlf2de
$ws3.Range("D3").Value = $handoffDe
$ws3.Range("H3").Value = "Include"

$ws3.Range("A4").Value = ".localization-config"
$ws3.Range("B4").Value = "Not to be localized"
$ws3.Range("C4").ClearContents()
$ws3.Range("D4").Value = $epoch
$ws3.Range("H4").Value = "Ignored"

$ws3.Range("I2").ClearContents()
$ws3.Range("I3").ClearContents()

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $urlMd1, "", "", $md1)
$ws3.Hyperlinks.Add($ws3.Range("C2"), $urlXlf1De, "", "", $xlf1de)
$ws3.Hyperlinks.Add($ws3.Range("A3"), $urlMd2, "", "", $md2)
$ws3.Hyperlinks.Add($ws3.Range("C3"), $urlXlf2De, "", "", $xlf2de)
$ws3.Hyperlinks.Add($ws3.Range("A4"), $urlConfig, "", "", ".localization-config")

$wb.Save()
